# Unhide / show-detail the outline group that was collapsed over rows 3-22
# (the group's summary/header row 13 was itself collapsed too).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 3; $r -le 22; $r++) {
    $ws.Rows.Item($r).Hidden = $false
}

# Row 13 was the collapsed outline-summary row for the nested group - un-collapse it.
$ws.Rows.Item(13).Hidden = $false

# Row 23 held a duplicate "MR2.9 - Location Selection" entry (incl. "Barleen" as the
# Executing Party) that needs to be wiped back to an empty (but still styled) row.
$ws.Range("A23:AA23").ClearContents()

# Re-collapse/select row 23 as the active selection, matching the edited file's
# recorded UI state (row 23 selected after the clear).
$ws.Rows.Item(23).Select() | Out-Null
